$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.294.35'
$ws.Range('E2').Value = '  -1.70%  '
$ws.Range('D3').Value = '2.577.52'
$ws.Range('E3').Value = '  -2.74%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = "'588.77"
$ws.Range('E5').Value = '  -3.43%  '
$ws.Range('D6').Value = "'150.72"
$ws.Range('E6').Value = '  +1.47%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('E8').Value = '  -1.02%  '
$ws.Range('E9').Value = '  +0.19%  '
$ws.Range('D10').Value = "'5.71"
$ws.Range('E10').Value = '  +1.80%  '
$ws.Range('E11').Value = '  -0.77%  '
$ws.Range('E12').Value = '  -0.58%  '
$ws.Range('D13').Value = "'27.60"
$ws.Range('E13').Value = '  -0.40%  '
$ws.Range('D14').Value = '3.038.96'
$ws.Range('E14').Value = '  -2.86%  '
$ws.Range('D15').Value = '63.124.60'
$ws.Range('E15').Value = '  -1.73%  '
$ws.Range('D16').Value = "'0.0000156"
$ws.Range('E16').Value = '  +5.28%  '
$ws.Range('D17').Value = '2.589.25'
$ws.Range('E17').Value = '  -2.52%  '
$ws.Range('D18').Value = "'12.19"
$ws.Range('E18').Value = '  +2.67%  '
$ws.Range('D19').Value = "'4.73"
$ws.Range('E19').Value = '  +2.77%  '
$ws.Range('D20').Value = "'345.44"
$ws.Range('E20').Value = '  -0.50%  '
$ws.Range('D21').Value = "'6.84"
$ws.Range('E21').Value = '  -1.13%  '
$ws.Range('D23').Value = "'67.11"
$ws.Range('E23').Value = '  +0.51%  '
$ws.Range('E24').Value = '  +1.10%  '
$ws.Range('D25').Value = "'9.12"
$ws.Range('E25').Value = '  -3.10%  '
$ws.Range('E26').Value = '  -2.90%  '
$ws.Range('D27').Value = "'558.14"
$ws.Range('E27').Value = '  +0.00%  '
$ws.Range('E28').Value = '  -1.56%  '
$ws.Range('D29').Value = "'0.162"
$ws.Range('E29').Value = '  +0.80%  '
$ws.Range('E30').Value = '  +0.27%  '
$ws.Range('E31').Value = '  -2.15%  '
$ws.Range('D32').Value = '0.0₃0857'
$ws.Range('E32').Value = '  -0.35%  '
$ws.Range('E33').Value = '  -1.39%  '
$ws.Range('D34').Value = "'5.21"
$ws.Range('E34').Value = '  -1.89%  '
$ws.Range('D35').Value = "'166.17"
$ws.Range('E35').Value = '  -1.98%  '
$ws.Range('E36').Value = '  +1.15%  '
$ws.Range('E37').Value = '  -0.27%  '
$ws.Range('D38').Value = "'19.51"
$ws.Range('E38').Value = '  +0.72%  '
$ws.Range('E39').Value = '  -2.09%  '
$ws.Range('E40').Value = '  -0.07%  '
$ws.Range('D41').Value = "'165.27"
$ws.Range('E41').Value = '  -0.20%  '
$ws.Range('E42').Value = '  -1.69%  '
$ws.Range('E43').Value = '  +3.63%  '
$ws.Range('D44').Value = "'22.87"
$ws.Range('E44').Value = '  +3.24%  '
$ws.Range('D45').Value = "'0.0585"
$ws.Range('E45').Value = '  +2.24%  '
$ws.Range('E46').Value = '  +5.80%  '
$ws.Range('E47').Value = '  -0.37%  '
$ws.Range('E48').Value = '  +1.16%  '
$ws.Range('D49').Value = "'0.0961"
$ws.Range('E49').Value = '  -0.03%  '
$ws.Range('D50').Value = "'19.09"
$ws.Range('E50').Value = '  +0.51%  '
$ws.Range('D51').Value = '0.0₆0234'
$ws.Range('E51').Value = '  +19.32%  '

# Remove the quote-prefix text style Excel applied so these cells match
# the original (unstyled) inline-string cells exactly.
$ws.Range('D5').ClearFormats()
$ws.Range('D6').ClearFormats()
$ws.Range('D10').ClearFormats()
$ws.Range('D13').ClearFormats()
$ws.Range('D16').ClearFormats()
$ws.Range('D18').ClearFormats()
$ws.Range('D19').ClearFormats()
$ws.Range('D20').ClearFormats()
$ws.Range('D21').ClearFormats()
$ws.Range('D23').ClearFormats()
$ws.Range('D25').ClearFormats()
$ws.Range('D27').ClearFormats()
$ws.Range('D29').ClearFormats()
$ws.Range('D34').ClearFormats()
$ws.Range('D35').ClearFormats()
$ws.Range('D38').ClearFormats()
$ws.Range('D41').ClearFormats()
$ws.Range('D44').ClearFormats()
$ws.Range('D45').ClearFormats()
$ws.Range('D49').ClearFormats()
$ws.Range('D50').ClearFormats()
